$wb = $excel.ActiveWorkbook

# Helper: force a value to be stored as literal text, bypassing Excel's
# automatic "looks like a date" inference for ISO-style yyyy-MM-dd strings.
# Writes the text into a scratch cell formatted as Text ("@"), copies just the
# VALUE (not the format) onto the destination, then wipes the scratch cell so
# it leaves no trace in the sheet's used range / dimension.
function Set-TextValue {
    param($ws, $cellRef, $text)
    $scratch = $ws.Range("ZZ100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
    $ws.Application.CutCopyMode = $false
}

# ---- Sheet 4: 保險 (insurance) ----
$ws4 = $wb.Worksheets.Item(4)

# Header row (row 1) -- copy the existing bold/bordered header style onto the newly added cells
$ws4.Range("B1").Copy()
$ws4.Range("E1:K1").PasteSpecial(-4122)
$ws4.Application.CutCopyMode = $false

$ws4.Range("B1").Value = "company"
$ws4.Range("C1").Value = "name"
$ws4.Range("D1").Value = "owner"
$ws4.Range("E1").Value = "property_category"
$ws4.Range("F1").Value = "category"
$ws4.Range("G1").Value = "date"
$ws4.Range("H1").Value = "legislator_name"
$ws4.Range("I1").Value = "legislator_id"
$ws4.Range("J1").Value = "source_file"
$ws4.Range("K1").Value = "index"

# Data rows (row 2 and row 3)
$ws4.Range("B2").Value = "台灣人壽"
$ws4.Range("C2").Value = "台灣人壽歲歲長泰還本终身險"
$ws4.Range("D2").Value = "李永得"
$ws4.Range("E2").Value = "insurance"
$ws4.Range("F2").Value = "normal"
Set-TextValue $ws4 "G2" "2011-11-25"
$ws4.Range("H2").Value = "邱議瑩"
$ws4.Range("I2").Value = 913
$ws4.Range("J2").Value = "tmp5431"
$ws4.Range("K2").Value = 49

$ws4.Range("B3").Value = "台灣人壽"
$ws4.Range("C3").Value = "台灣人壽新祥和定期壽險"
$ws4.Range("D3").Value = "李永得"
$ws4.Range("E3").Value = "insurance"
$ws4.Range("F3").Value = "normal"
Set-TextValue $ws4 "G3" "2011-11-25"
$ws4.Range("H3").Value = "邱議瑩"
$ws4.Range("I3").Value = 913
$ws4.Range("J3").Value = "tmp5431"
$ws4.Range("K3").Value = 50

# ---- Sheet 5: 事業投資 (business investment) ----
$ws5 = $wb.Worksheets.Item(5)

# Header row (row 1) -- copy the existing bold/bordered header style onto the newly added cells
$ws5.Range("B1").Copy()
$ws5.Range("H1:N1").PasteSpecial(-4122)
$ws5.Application.CutCopyMode = $false

$ws5.Range("B1").Value = "owner"
$ws5.Range("C1").Value = "company"
$ws5.Range("D1").Value = "address"
$ws5.Range("E1").Value = "total"
$ws5.Range("F1").Value = "register_date"
$ws5.Range("G1").Value = "register_reason"
$ws5.Range("H1").Value = "property_category"
$ws5.Range("I1").Value = "category"
$ws5.Range("J1").Value = "date"
$ws5.Range("K1").Value = "legislator_name"
$ws5.Range("L1").Value = "legislator_id"
$ws5.Range("M1").Value = "source_file"
$ws5.Range("N1").Value = "index"

# Data rows (rows 2-4)
$ws5.Range("B2").Value = "李永得"
$ws5.Range("C2").Value = "玉山社事業股份有限公司"
$ws5.Range("D2").Value = "臺北市大安區仁愛路四段145號3樓之2"
$ws5.Range("E2").Value = 1000000
$ws5.Range("F2").Value = "84年07月08日"
$ws5.Range("G2").Value = "發起設立"
$ws5.Range("H2").Value = "investment"
$ws5.Range("I2").Value = "normal"
Set-TextValue $ws5 "J2" "2011-11-25"
$ws5.Range("K2").Value = "邱議瑩"
$ws5.Range("L2").Value = 913
$ws5.Range("M2").Value = "tmp5431"
$ws5.Range("N2").Value = 55

$ws5.Range("B3").Value = "李永得"
$ws5.Range("C3").Value = "淡海股份有限公司"
$ws5.Range("D3").Value = "新北市淡水區真理街17號"
$ws5.Range("E3").Value = 280800
$ws5.Range("F3").Value = "97年02月27日"
$ws5.Range("G3").Value = "91.07.25發起"
$ws5.Range("H3").Value = "investment"
$ws5.Range("I3").Value = "normal"
Set-TextValue $ws5 "J3" "2011-11-25"
$ws5.Range("K3").Value = "邱議瑩"
$ws5.Range("L3").Value = 913
$ws5.Range("M3").Value = "tmp5431"
$ws5.Range("N3").Value = 56

$ws5.Range("B4").Value = "李永得"
$ws5.Range("C4").Value = "先驅媒體社會企業股份有限公司"
$ws5.Range("D4").Value = "臺北市中正區仁愛路二段98號7樓"
$ws5.Range("E4").Value = 500000
$ws5.Range("F4").Value = "98年07月24日"
$ws5.Range("G4").Value = "發起"
$ws5.Range("H4").Value = "investment"
$ws5.Range("I4").Value = "normal"
Set-TextValue $ws5 "J4" "2011-11-25"
$ws5.Range("K4").Value = "邱議瑩"
$ws5.Range("L4").Value = 913
$ws5.Range("M4").Value = "tmp5431"
$ws5.Range("N4").Value = 57

